$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: in-place character edits so the shared-string run structure is preserved as much as possible ---
# A8 = "Volume 31   Number  1" -> "...Number  2" (char 21 is the "1")
$ws.Range("A8").Characters(21, 1).Text = "2"
# C9 = "Report Covering the Week  1/1/2024  Through  1/7/2024" -> week shifted forward 7 days
# edit right-most substring first so the left offset (27) stays valid
$ws.Range("C9").Characters(46, 8).Text = "1/14/2024"
$ws.Range("C9").Characters(27, 8).Text = "1/8/2024"

# --- Value-only updates: style & data type are unchanged, just refresh the figure ---
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -61.904761904761
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 12
$ws.Range("K16").Value = -58.333333333333
$ws.Range("L16").Value = -28.571428571428
$ws.Range("M16").Value = 25
$ws.Range("N16").Value = -86.111111111111
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 13
$ws.Range("H17").Value = 18.181818181818
$ws.Range("I17").Value = 6
$ws.Range("J17").Value = 7
$ws.Range("K17").Value = -14.285714285714
$ws.Range("L17").Value = 50
$ws.Range("M17").Value = 500
$ws.Range("N17").Value = -50
$ws.Range("C18").Value = 7
$ws.Range("E18").Value = -22.222222222222
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 33
$ws.Range("H18").Value = -27.272727272727
$ws.Range("I18").Value = 9
$ws.Range("J18").Value = 18
$ws.Range("K18").Value = -50
$ws.Range("L18").Value = -47.058823529411
$ws.Range("M18").Value = -10
$ws.Range("N18").Value = -66.666666666666
$ws.Range("C19").Value = 26
$ws.Range("D19").Value = 28
$ws.Range("E19").Value = -7.142857142857
$ws.Range("F19").Value = 82
$ws.Range("G19").Value = 89
$ws.Range("H19").Value = -7.865168539325
$ws.Range("I19").Value = 43
$ws.Range("J19").Value = 54
$ws.Range("K19").Value = -20.370370370370
$ws.Range("L19").Value = 13.157894736842
$ws.Range("M19").Value = 13.157894736842
$ws.Range("N19").Value = -60.185185185185
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = 150
$ws.Range("N20").Value = -93.103448275862
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 46
$ws.Range("E21").Value = -13.043478260869
$ws.Range("F21").Value = 133
$ws.Range("G21").Value = 157
$ws.Range("H21").Value = -15.286624203821
$ws.Range("I21").Value = 65
$ws.Range("J21").Value = 93
$ws.Range("K21").Value = -30.107526881720
$ws.Range("L21").Value = -9.722222222222
$ws.Range("M21").Value = 20.370370370370
$ws.Range("N21").Value = -69.339622641509
$ws.Range("C22").Value = 1
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 25
$ws.Range("I22").Value = 4
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -42.857142857142
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -21.428571428571
$ws.Range("F24").Value = 121
$ws.Range("G24").Value = 127
$ws.Range("H24").Value = -4.724409448818
$ws.Range("I24").Value = 51
$ws.Range("J24").Value = 61
$ws.Range("K24").Value = -16.393442622950
$ws.Range("L24").Value = -17.741935483871
$ws.Range("M24").Value = 18.604651162790
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 66.666666666666
$ws.Range("F25").Value = 34
$ws.Range("H25").Value = 41.666666666666
$ws.Range("I25").Value = 18
$ws.Range("J25").Value = 16
$ws.Range("K25").Value = 12.5
$ws.Range("L25").Value = 200
$ws.Range("M25").Value = 200
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -16.666666666666
$ws.Range("J27").Value = 4
$ws.Range("K27").Value = -50
$ws.Range("J42").Value = 48
$ws.Range("K42").Value = -70.909090909090
$ws.Range("L42").Value = -82.022471910112
$ws.Range("M42").Value = -93.790426908150
$ws.Range("N42").Value = -95.604395604395
$ws.Range("J43").Value = 2067
$ws.Range("K43").Value = -2.5
$ws.Range("L43").Value = -29.837067209776
$ws.Range("M43").Value = -61.277632071937
$ws.Range("N43").Value = -73.628476652207

# --- Cells whose data type flips between number and text (zero/insufficient-data placeholders) ---
# Each one: write the new value (prefixing text with an apostrophe so it is not re-parsed as a number),
# then copy *just the format* from a donor cell that already carries the destination style, so the
# cell lands on the exact same style index the workbook already uses for that look (s=14 text / s=15 int / s=16 pct).
$ws.Range("D15").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").Value = "'***.*"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("L15").Value = -100
$ws.Range("H15").Copy() | Out-Null
$ws.Range("L15").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("D20").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = "'***.*"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("L20").Value = -50
$ws.Range("H15").Copy() | Out-Null
$ws.Range("L20").PasteSpecial(-4122) | Out-Null
$ws.Range("M20").Value = 100
$ws.Range("H15").Copy() | Out-Null
$ws.Range("M20").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = 0
$ws.Range("H15").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("J22").Value = 1
$ws.Range("F15").Copy() | Out-Null
$ws.Range("J22").PasteSpecial(-4122) | Out-Null
$ws.Range("K22").Value = 300
$ws.Range("H15").Copy() | Out-Null
$ws.Range("K22").PasteSpecial(-4122) | Out-Null
$ws.Range("D26").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Value = "'***.*"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("L26").Value = -100
$ws.Range("H15").Copy() | Out-Null
$ws.Range("L26").PasteSpecial(-4122) | Out-Null
$ws.Range("C27").Value = 2
$ws.Range("F15").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("I27").Value = 2
$ws.Range("F15").Copy() | Out-Null
$ws.Range("I27").PasteSpecial(-4122) | Out-Null
$ws.Range("L27").Value = 100
$ws.Range("H15").Copy() | Out-Null
$ws.Range("L27").PasteSpecial(-4122) | Out-Null
